# Auto-generated edit script applying numeric corrections to the
# Leve profit-tracking columns (H-N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR
# sheets, per the scheduled pricing-refresh run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 19: Unbreak My Heart
$ws.Range("H19").Value = 8602.385
$ws.Range("I19").Value = 50455.5
$ws.Range("J19").Value = 992.7273
$ws.Range("K19").Value = 50455.5
$ws.Range("L19").Value = 992.7273
$ws.Range("M19").Value = -50280.5
$ws.Range("N19").Value = -1342.7273

# ALC row 130: Technically Still Magic
$ws.Range("H130").Value = 12153.685
$ws.Range("J130").Value = 12153.685
$ws.Range("L130").Value = 12153.685
$ws.Range("N130").Value = -22193.685

# ALC row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 1641.0476
$ws.Range("I137").Value = 1741.375
$ws.Range("K137").Value = 5224.125
$ws.Range("M137").Value = -2674.125

$ws = $wb.Worksheets.Item("ARM")
# ARM row 45: Hollow Hallmarks
$ws.Range("H45").Value = 2568.75
$ws.Range("I45").Value = 916.6667
$ws.Range("J45").Value = 3560
$ws.Range("K45").Value = 916.6667
$ws.Range("L45").Value = 3560
$ws.Range("M45").Value = -539.6667
$ws.Range("N45").Value = -4314

# ARM row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 3472.5789
$ws.Range("I61").Value = 2638.4546
$ws.Range("K61").Value = 2638.4546
$ws.Range("M61").Value = -2426.4546

# ARM row 74: As the Bolt Flies
$ws.Range("H74").Value = 1301.3438
$ws.Range("I74").Value = 978.7222
$ws.Range("J74").Value = 1716.1428
$ws.Range("K74").Value = 978.7222
$ws.Range("L74").Value = 1716.1428
$ws.Range("M74").Value = -104.7222
$ws.Range("N74").Value = -3464.1428

# ARM row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value = 1301.3438
$ws.Range("I77").Value = 978.7222
$ws.Range("J77").Value = 1716.1428
$ws.Range("K77").Value = 4893.611
$ws.Range("L77").Value = 8580.714
$ws.Range("M77").Value = -525.6109999999999
$ws.Range("N77").Value = -17316.714

# ARM row 97: Ore for Me
$ws.Range("H97").Value = 584.4706
$ws.Range("I97").Value = 410.68967
$ws.Range("J97").Value = 1592.4
$ws.Range("K97").Value = 410.68967
$ws.Range("L97").Value = 1592.4
$ws.Range("M97").Value = 85.31033000000002
$ws.Range("N97").Value = -2584.4

# ARM row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 4246.8335
$ws.Range("I132").Value = 2513.1
$ws.Range("K132").Value = 7539.299999999999
$ws.Range("M132").Value = -5009.299999999999

# ARM row 136: Metal with Mettle
$ws.Range("H136").Value = 3472.5789
$ws.Range("I136").Value = 2638.4546
$ws.Range("K136").Value = 7915.3638
$ws.Range("M136").Value = -5365.3638

$ws = $wb.Worksheets.Item("BSM")
# BSM row 86: Through Thick and Thin
$ws.Range("H86").Value = 1873.2609
$ws.Range("I86").Value = 1891.1538
$ws.Range("J86").Value = 1850
$ws.Range("K86").Value = 1891.1538
$ws.Range("L86").Value = 1850
$ws.Range("M86").Value = -768.1538
$ws.Range("N86").Value = -4096

# BSM row 89: Piercing Eyes Deserve Piercing Shafts (L)
$ws.Range("H89").Value = 1873.2609
$ws.Range("I89").Value = 1891.1538
$ws.Range("J89").Value = 1850
$ws.Range("K89").Value = 9455.769
$ws.Range("L89").Value = 9250
$ws.Range("M89").Value = -3839.769
$ws.Range("N89").Value = -20482

$ws = $wb.Worksheets.Item("CRP")
# CRP row 31: Wall Not Found
$ws.Range("H31").Value = 2124.75
$ws.Range("I31").Value = 729.4054
$ws.Range("J31").Value = 2944.238
$ws.Range("K31").Value = 729.4054
$ws.Range("L31").Value = 2944.238
$ws.Range("M31").Value = -434.4054
$ws.Range("N31").Value = -3534.238

# CRP row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 2124.75
$ws.Range("I34").Value = 729.4054
$ws.Range("J34").Value = 2944.238
$ws.Range("K34").Value = 729.4054
$ws.Range("L34").Value = 2944.238
$ws.Range("M34").Value = -527.4054
$ws.Range("N34").Value = -3348.238

# CRP row 58: You Do the Heavy Lifting
$ws.Range("H58").Value = 848.7778
$ws.Range("I58").Value = 819.26666
$ws.Range("J58").Value = 996.3333
$ws.Range("K58").Value = 819.26666
$ws.Range("L58").Value = 996.3333
$ws.Range("M58").Value = -616.26666
$ws.Range("N58").Value = -1402.3333

# CRP row 134: Wood You Be Quiet
$ws.Range("H134").Value = 3638.925
$ws.Range("I134").Value = 4534.64
$ws.Range("K134").Value = 13603.92
$ws.Range("M134").Value = -11068.92

# CRP row 136: Turali Quality
$ws.Range("H136").Value = 848.7778
$ws.Range("I136").Value = 819.26666
$ws.Range("J136").Value = 996.3333
$ws.Range("K136").Value = 2457.79998
$ws.Range("L136").Value = 2988.9999
$ws.Range("M136").Value = 92.20002000000022
$ws.Range("N136").Value = -8088.9999

$ws = $wb.Worksheets.Item("CUL")
# CUL row 4: In Hot Water
$ws.Range("H4").Value = 742.6316
$ws.Range("I4").Value = 177.6923
$ws.Range("J4").Value = 1966.6666
$ws.Range("K4").Value = 533.0769
$ws.Range("L4").Value = 5899.9998
$ws.Range("M4").Value = -421.0769
$ws.Range("N4").Value = -6123.9998

# CUL row 5: What a Sap
$ws.Range("H5").Value = 338481.66
$ws.Range("I5").Value = 739.2
$ws.Range("J5").Value = 418896.53
$ws.Range("K5").Value = 2217.6
$ws.Range("L5").Value = 1256689.59
$ws.Range("M5").Value = -2105.6
$ws.Range("N5").Value = -1256913.59

# CUL row 68: Such a Butter Face
$ws.Range("H68").Value = 2138200.8
$ws.Range("I68").Value = 6740784
$ws.Range("J68").Value = 1287.0714
$ws.Range("K68").Value = 20222352
$ws.Range("L68").Value = 3861.2142
$ws.Range("M68").Value = -20221541
$ws.Range("N68").Value = -5483.2142

# CUL row 71: No Margarine of Error (L)
$ws.Range("H71").Value = 2138200.8
$ws.Range("I71").Value = 6740784
$ws.Range("J71").Value = 1287.0714
$ws.Range("K71").Value = 60667056
$ws.Range("L71").Value = 11583.6426
$ws.Range("M71").Value = -60663000
$ws.Range("N71").Value = -19695.6426

# CUL row 113: Can't Eat Just One
$ws.Range("H113").Value = 2286.6875
$ws.Range("I113").Value = 3609.6667
$ws.Range("J113").Value = 1492.9
$ws.Range("K113").Value = 10829.0001
$ws.Range("L113").Value = 4478.700000000001
$ws.Range("M113").Value = -8659.000100000001
$ws.Range("N113").Value = -8818.700000000001

# CUL row 134: Don't Knock It Till You've Tried It
$ws.Range("H134").Value = 39891.69
$ws.Range("I134").Value = 67645.60000000001
$ws.Range("J134").Value = 2045.4546
$ws.Range("K134").Value = 202936.8
$ws.Range("L134").Value = 6136.3638
$ws.Range("M134").Value = -197866.8
$ws.Range("N134").Value = -16276.3638

# CUL row 135: Not-so-secret Ingredient
$ws.Range("H135").Value = 338481.66
$ws.Range("I135").Value = 739.2
$ws.Range("J135").Value = 418896.53
$ws.Range("K135").Value = 6652.8
$ws.Range("L135").Value = 3770068.77
$ws.Range("M135").Value = -4117.8
$ws.Range("N135").Value = -3775138.77

# CUL row 139: Najoothie
$ws.Range("H139").Value = 35757.9
$ws.Range("I139").Value = 59598.766
$ws.Range("J139").Value = 1983.3334
$ws.Range("K139").Value = 178796.298
$ws.Range("L139").Value = 5950.0002
$ws.Range("M139").Value = -173656.298
$ws.Range("N139").Value = -16230.0002

# CUL row 140: Sweet, Sweet Bean Juice
$ws.Range("H140").Value = 34546.234
$ws.Range("I140").Value = 53394.05
$ws.Range("J140").Value = 1990.909
$ws.Range("K140").Value = 160182.15
$ws.Range("L140").Value = 5972.727000000001
$ws.Range("M140").Value = -155002.15
$ws.Range("N140").Value = -16332.727

$ws = $wb.Worksheets.Item("GSM")
# GSM row 132: On Board for Lar
$ws.Range("H132").Value = 5239
$ws.Range("I132").Value = 1966.3334
$ws.Range("J132").Value = 7202.6
$ws.Range("K132").Value = 5899.0002
$ws.Range("L132").Value = 21607.8
$ws.Range("M132").Value = -3369.0002
$ws.Range("N132").Value = -26667.8

$ws = $wb.Worksheets.Item("LTW")
# LTW row 2: Red in the Head
$ws.Range("H2").Value = 500000
$ws.Range("J2").Value = 500000
$ws.Range("L2").Value = 500000
$ws.Range("N2").Value = -500224

# LTW row 7: Tan Before the Ban
$ws.Range("H7").Value = 52743.4
$ws.Range("I7").Value = 65298
$ws.Range("J7").Value = 2525
$ws.Range("K7").Value = 65298
$ws.Range("L7").Value = 2525
$ws.Range("M7").Value = -65186
$ws.Range("N7").Value = -2749

# LTW row 22: Skin off Their Backs
$ws.Range("H22").Value = 1020.73334
$ws.Range("I22").Value = 979.4
$ws.Range("J22").Value = 1103.4
$ws.Range("K22").Value = 979.4
$ws.Range("L22").Value = 1103.4
$ws.Range("M22").Value = -684.4
$ws.Range("N22").Value = -1693.4

# LTW row 27: Fire and Hide
$ws.Range("H27").Value = 1020.73334
$ws.Range("I27").Value = 979.4
$ws.Range("J27").Value = 1103.4
$ws.Range("K27").Value = 979.4
$ws.Range("L27").Value = 1103.4
$ws.Range("M27").Value = -872.4
$ws.Range("N27").Value = -1317.4

# LTW row 46: Supply Side Logic
$ws.Range("H46").Value = 112755.555
$ws.Range("I46").Value = 144571.58
$ws.Range("J46").Value = 1399.5
$ws.Range("K46").Value = 144571.58
$ws.Range("L46").Value = 1399.5
$ws.Range("M46").Value = -144383.58
$ws.Range("N46").Value = -1775.5

# LTW row 100: Tiger in the Sack
$ws.Range("H100").Value = 3000
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()

# LTW row 122: Hell on Leather
$ws.Range("H122").Value = 13890636
$ws.Range("I122").Value = 18520102
$ws.Range("J122").Value = 2237.5
$ws.Range("K122").Value = 55560306
$ws.Range("L122").Value = 6712.5
$ws.Range("M122").Value = -55557856
$ws.Range("N122").Value = -11612.5

# LTW row 126: Battered Books
$ws.Range("H126").Value = 52743.4
$ws.Range("I126").Value = 65298
$ws.Range("J126").Value = 2525
$ws.Range("K126").Value = 195894
$ws.Range("L126").Value = 7575
$ws.Range("M126").Value = -193424
$ws.Range("N126").Value = -12515

$ws = $wb.Worksheets.Item("WVR")
# WVR row 2: The Unmentionables
$ws.Range("H2").Value = 68003
$ws.Range("J2").Value = 68003
$ws.Range("L2").Value = 68003
$ws.Range("N2").Value = -68227

# WVR row 92: Modest Beginnings
$ws.Range("H92").Value = 15000
$ws.Range("J92").Value = 15000
$ws.Range("L92").Value = 15000
$ws.Range("N92").Value = -19992

# WVR row 132: Comfy Cabins
$ws.Range("H132").Value = 2359.1333
$ws.Range("I132").Value = 1432.1765
$ws.Range("K132").Value = 4296.529500000001
$ws.Range("M132").Value = -1766.529500000001

